$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "DonacionesProyectos (2)"

# Select entire rows 1 through 6 (equivalent to A1:XFD6)
$ws.Rows("1:6").Select() | Out-Null
